$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 53) onto the new row 54
$ws.Range("A53:E53").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's values
$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 2.46481303148316
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 3.633434696013671
